# Clouducate.pptx — "AWSVPCB SCRIPTS" slide (slide 9) text revisions.
#
# 1) AWSVPCB.VPC.CREATE bullet: split the single explanatory run into three
#    runs so it reads "... This script runs VPC.REGISTER, creates your VPC,
#    Internet Gateway, ... AWS unique IDs " (trailing space added).
# 2) AWSVPCB.VPC.REGISTER bullet: merge the trailing "run automatically by
#    VPC.CREATE" run into the preceding explanatory run and append
#    " & VPC.DESTROY".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- AWSVPCB.VPC.CREATE (paragraph 4) ---------------------------------
$para4 = $tr.Paragraphs(4)

# Shorten the explanatory run down to its new opening clause.
$run2 = $para4.Runs(2)
$run2.Text = "– This script runs VPC.REGISTER"

# Add the ", creates " connector run right after it.
$run2 = $tr.Paragraphs(4).Runs(2)
$run2.InsertAfter(", creates ")

# Add the remainder of the original sentence (now with a trailing space)
# as its own run after the connector.
$run3 = $tr.Paragraphs(4).Runs(3)
$run3.InsertAfter("your VPC, Internet Gateway, route tables, subnets, security groups, S3 buckets, NAT instance, BASTION server, etc. and registers all AWS unique IDs ")

# --- AWSVPCB.VPC.REGISTER (paragraph 7) --------------------------------
$para7 = $tr.Paragraphs(7)

# Combine the explanatory run with the trailing "run automatically by
# VPC.CREATE" run into one run, then add the new "& VPC.DESTROY" mention.
$r2 = $para7.Runs(2)
$r3 = $para7.Runs(3)
$start = $r2.Start
$len = $r3.Start + $r3.Text.Length - $r2.Start
$combined = $tr.Characters($start, $len)
$combined.Text = "– May never be needed. It will find your AWS-VPCB VPC in AWS and register all its components for the other scripts to be able to work as expected – run automatically by VPC.CREATE & VPC.DESTROY"
